# Update gh-pages output generated at 456a3b4
# Updates the "F" column (attendee/follower counts) across the four sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 228
$ws1.Range("F4").Value  = 813
$ws1.Range("F6").Value  = 400
$ws1.Range("F7").Value  = 560
$ws1.Range("F8").Value  = 213
$ws1.Range("F9").Value  = 59
$ws1.Range("F11").Value = 128
$ws1.Range("F12").Value = 609
$ws1.Range("F13").Value = 78
$ws1.Range("F14").Value = 1760
$ws1.Range("F15").Value = 321
$ws1.Range("F16").Value = 2378
$ws1.Range("F17").Value = 295
$ws1.Range("F19").Value = 42
$ws1.Range("F20").Value = 129

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 221
$ws2.Range("F5").Value  = 13
$ws2.Range("F13").Value = 84

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5282
$ws3.Range("F3").Value = 305
$ws3.Range("F4").Value = 197

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5282
$ws4.Range("F4").Value  = 305
$ws4.Range("F6").Value  = 197
$ws4.Range("F7").Value  = 228
$ws4.Range("F8").Value  = 221
$ws4.Range("F10").Value = 13
$ws4.Range("F13").Value = 813
$ws4.Range("F17").Value = 400
$ws4.Range("F18").Value = 560
$ws4.Range("F19").Value = 213
$ws4.Range("F20").Value = 59
$ws4.Range("F23").Value = 128
$ws4.Range("F26").Value = 609
$ws4.Range("F27").Value = 78
$ws4.Range("F28").Value = 84
$ws4.Range("F29").Value = 1760
$ws4.Range("F30").Value = 321
$ws4.Range("F31").Value = 2379
$ws4.Range("F33").Value = 295
$ws4.Range("F35").Value = 42
$ws4.Range("F36").Value = 129

$wb.Save()
